# Part1-UML-Review.pptx - "updating various slides and fixing exception
# project setup for a single solution project"
#
# Slide 1 changes:
#  1. The "Rectangle 5" shape (the git-projects / quiz bullet list) loses the
#     "PracticeFilesAndExceptionsTooManyScoresSolution" bullet paragraph
#     (the exception project is now a single solution project) and is
#     resized/repositioned (shrinks a bit + moves down) to account for the
#     removed line.
#  2. The "TextBox 1" shape's attendance-password blank grows from 10
#     underscores to 11 underscores.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 3: "Rectangle 5" (the bulleted "git projects" / "Quiz" box) ---
$rect = $s.Shapes.Item(3)

# Remove the "PracticeFilesAndExceptionsTooManyScoresSolution" bullet line
# (paragraph 4 of its text: "today's projects:" / 3 bullets / "quiz:" / bullet)
# Note: Paragraphs().Text includes the trailing paragraph-mark (CR) for every
# paragraph except the very last one in the shape, so trim before comparing.
$paragraphs = $rect.TextFrame.TextRange.Paragraphs()
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $rect.TextFrame.TextRange.Paragraphs($i)
    if ($para.Text.Trim() -eq "PracticeFilesAndExceptionsTooManyScoresSolution") {
        $para.Delete()
        break
    }
}

# Reposition/resize the box now that it holds one fewer line of text.
$rect.Top = 367.0286
$rect.Height = 159.7124

# --- Shape 4: "TextBox 1" (attendance password blank) ---
$box = $s.Shapes.Item(4)
$boxParagraphs = $box.TextFrame.TextRange.Paragraphs()
for ($i = 1; $i -le $boxParagraphs.Count; $i++) {
    $para = $box.TextFrame.TextRange.Paragraphs($i)
    if ($para.Text.Trim() -eq "__________") {
        $para.Runs(1).Text = "___________"
        break
    }
}
